$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "42.470.84"
$ws.Cells.Item(2,5).Value = "  -1.28%  "
$ws.Cells.Item(3,4).Value = "2.227.02"
$ws.Cells.Item(3,5).Value = "  -0.60%  "
$ws.Cells.Item(4,5).Value = "  -0.05%  "
$ws.Cells.Item(5,4).Value = "'111.99"
$ws.Cells.Item(5,5).Value = "  -3.79%  "
$ws.Cells.Item(6,4).Value = "'291.42"
$ws.Cells.Item(6,5).Value = "  +9.31%  "
$ws.Cells.Item(7,5).Value = "  -0.80%  "
$ws.Cells.Item(8,5).Value = "  -0.30%  "
$ws.Cells.Item(9,4).Value = "'0.601"
$ws.Cells.Item(9,5).Value = "  -1.51%  "
$ws.Cells.Item(10,4).Value = "'43.88"
$ws.Cells.Item(10,5).Value = "  -6.06%  "
$ws.Cells.Item(11,4).Value = "'0.0913"
$ws.Cells.Item(11,5).Value = "  -1.93%  "
$ws.Cells.Item(12,4).Value = "'54.44"
$ws.Cells.Item(12,5).Value = "  +1.14%  "
$ws.Cells.Item(13,4).Value = "'8.66"
$ws.Cells.Item(13,5).Value = "  -5.70%  "
$ws.Cells.Item(14,4).Value = "'1.03"
$ws.Cells.Item(14,5).Value = "  +16.79%  "
$ws.Cells.Item(15,5).Value = "  -1.73%  "
$ws.Cells.Item(16,4).Value = "'14.95"
$ws.Cells.Item(16,5).Value = "  -2.78%  "
$ws.Cells.Item(17,4).Value = "2.560.33"
$ws.Cells.Item(17,5).Value = "  -0.70%  "
$ws.Cells.Item(18,4).Value = "2.222.17"
$ws.Cells.Item(18,5).Value = "  -1.91%  "
$ws.Cells.Item(19,4).Value = "42.473.50"
$ws.Cells.Item(19,5).Value = "  -1.42%  "
$ws.Cells.Item(20,2).Value = "Uniswap"
$ws.Cells.Item(20,3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(20,4).Value = "'7.17"
$ws.Cells.Item(20,5).Value = "  +6.29%  "
$ws.Cells.Item(21,2).Value = "ShibaInu"
$ws.Cells.Item(21,3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(21,4).Value = "'0.0000105"
$ws.Cells.Item(21,5).Value = "  -2.19%  "
$ws.Cells.Item(22,4).Value = "'73.54"
$ws.Cells.Item(22,5).Value = "  +2.69%  "
$ws.Cells.Item(23,4).Value = "'3.33"
$ws.Cells.Item(23,5).Value = "  +14.35%  "
$ws.Cells.Item(24,4).Value = "'2.39"
$ws.Cells.Item(24,5).Value = "  +0.64%  "
$ws.Cells.Item(25,4).Value = "'235.73"
$ws.Cells.Item(25,5).Value = "  +1.73%  "
$ws.Cells.Item(26,4).Value = "'8.92"
$ws.Cells.Item(26,5).Value = "  -5.67%  "
$ws.Cells.Item(27,4).Value = "'1.00"
$ws.Cells.Item(27,5).Value = "  -1.81%  "
$ws.Cells.Item(28,5).Value = "  -7.41%  "
$ws.Cells.Item(29,5).Value = "  -1.70%  "
$ws.Cells.Item(30,4).Value = "'37.70"
$ws.Cells.Item(30,5).Value = "  -8.45%  "
$ws.Cells.Item(31,4).Value = "'173.72"
$ws.Cells.Item(32,4).Value = "'3.14"
$ws.Cells.Item(32,5).Value = "  -4.85%  "
$ws.Cells.Item(33,4).Value = "'21.28"
$ws.Cells.Item(33,5).Value = "  +0.13%  "
$ws.Cells.Item(34,4).Value = "'0.0878"
$ws.Cells.Item(35,4).Value = "'5.62"
$ws.Cells.Item(35,5).Value = "  +0.36%  "
$ws.Cells.Item(36,4).Value = "'5.02"
$ws.Cells.Item(36,5).Value = "  +7.08%  "
$ws.Cells.Item(37,2).Value = "NEARProtocol"
$ws.Cells.Item(37,3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(37,4).Value = "'4.22"
$ws.Cells.Item(37,5).Value = "  -2.20%  "
$ws.Cells.Item(38,2).Value = "Stellar"
$ws.Cells.Item(38,3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(38,4).Value = "'0.126"
$ws.Cells.Item(38,5).Value = "  -1.74%  "
$ws.Cells.Item(39,5).Value = "  +0.43%  "
$ws.Cells.Item(40,4).Value = "'0.105"
$ws.Cells.Item(40,5).Value = "  -1.59%  "
$ws.Cells.Item(41,4).Value = "'2.39"
$ws.Cells.Item(41,5).Value = "  -4.59%  "
$ws.Cells.Item(42,4).Value = "'71.87"
$ws.Cells.Item(42,5).Value = "  +0.79%  "
$ws.Cells.Item(43,4).Value = "'0.231"
$ws.Cells.Item(43,5).Value = "  -1.62%  "
$ws.Cells.Item(44,5).Value = "  -0.10%  "
$ws.Cells.Item(45,4).Value = "'12.35"
$ws.Cells.Item(45,5).Value = "  -8.00%  "
$ws.Cells.Item(46,4).Value = "'1.31"
$ws.Cells.Item(46,5).Value = "  -2.16%  "
$ws.Cells.Item(47,4).Value = "'5.36"
$ws.Cells.Item(48,4).Value = "'1.28"
$ws.Cells.Item(48,5).Value = "  +2.91%  "
$ws.Cells.Item(49,2).Value = "Stacks"
$ws.Cells.Item(49,3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(49,4).Value = "'1.64"
$ws.Cells.Item(49,5).Value = "  +4.11%  "
$ws.Cells.Item(50,4).Value = "'8.42"
$ws.Cells.Item(50,5).Value = "  -0.14%  "
$ws.Cells.Item(51,2).Value = "Aave"
$ws.Cells.Item(51,3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(51,4).Value = "'101.08"
$ws.Cells.Item(51,5).Value = "  +0.55%  "
